# Updated cryptos list on Fri Dec 29 06:46:58 UTC 2023 with GitHub Actions
#
# Refresh the Price (col D) and Volume(1h) (col E) figures for each coin row,
# and swap row 51 from Cronos to TrustWalletToken.
#
# Note: several new Price values (e.g. "8.50", "1.00", "16.00") look like
# plain numbers to Excel's auto-detection and would otherwise be coerced to
# numeric cells (dropping the trailing zero / formatting). Those are written
# with a leading apostrophe so they are stored as literal text, matching the
# source data which keeps prices as text strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.665.03"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").Value = "2.359.19"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'319.03"
$ws.Range("E5").Value = "  -1.52%  "

$ws.Range("D6").Value = "'108.84"
$ws.Range("E6").Value = "  +3.16%  "

$ws.Range("E7").Value = "  -1.35%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -4.18%  "

$ws.Range("D10").Value = "'41.73"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").Value = "'8.50"
$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").Value = "'1.00"
$ws.Range("E13").Value = "  -1.99%  "

$ws.Range("E14").Value = "  +0.15%  "

$ws.Range("D15").Value = "'16.00"
$ws.Range("E15").Value = "  -6.97%  "

$ws.Range("D16").Value = "2.712.34"
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("D17").Value = "2.316.31"
$ws.Range("E17").Value = "  -3.71%  "

$ws.Range("D18").Value = "42.781.85"
$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("D19").Value = "'7.76"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").Value = "'0.0000107"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("D21").Value = "'76.70"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").Value = "'3.63"
$ws.Range("E22").Value = "  +5.92%  "

$ws.Range("D23").Value = "'257.36"
$ws.Range("E23").Value = "  -6.26%  "

$ws.Range("D24").Value = "'2.32"
$ws.Range("E24").Value = "  -3.55%  "

$ws.Range("D25").Value = "'9.40"
$ws.Range("E25").Value = "  -3.40%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "'11.46"
$ws.Range("E27").Value = "  -2.44%  "

$ws.Range("D28").Value = "'22.94"
$ws.Range("E28").Value = "  +0.52%  "

$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").Value = "'175.08"
$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("D31").Value = "'36.85"
$ws.Range("E31").Value = "  -2.72%  "

$ws.Range("D32").Value = "'0.0892"
$ws.Range("E32").Value = "  -3.55%  "

$ws.Range("E33").Value = "  +3.73%  "

$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  -8.62%  "

$ws.Range("E35").Value = "  +19.12%  "

$ws.Range("E36").Value = "  -1.15%  "

$ws.Range("D37").Value = "'4.63"
$ws.Range("E37").Value = "  -4.77%  "

$ws.Range("D38").Value = "'0.0363"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").Value = "'3.83"
$ws.Range("E39").Value = "  -6.77%  "

$ws.Range("E40").Value = "  -5.47%  "

$ws.Range("E41").Value = "  +3.18%  "

$ws.Range("D42").Value = "'72.08"
$ws.Range("E42").Value = "  +4.25%  "

$ws.Range("E43").Value = "  -6.77%  "

$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "'113.84"
$ws.Range("E45").Value = "  -8.45%  "

$ws.Range("D46").Value = "'12.02"
$ws.Range("E46").Value = "  -3.73%  "

$ws.Range("D47").Value = "'5.51"
$ws.Range("E47").Value = "  -1.75%  "

$ws.Range("D48").Value = "'9.14"
$ws.Range("E48").Value = "  -3.96%  "

$ws.Range("D49").Value = "'83.92"
$ws.Range("E49").Value = "  -9.46%  "

$ws.Range("D50").Value = "'73.77"
$ws.Range("E50").Value = "  +0.71%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.26"
$ws.Range("E51").Value = "  -3.81%  "
